$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL) ---
$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"

# --- Numeric-looking text columns (Price / Volume) ---
# Force text storage via NumberFormat "@" so values like "41.60", "0.3500"
# and "0.06%" are preserved verbatim (not coerced into numbers/percentages),
# then restore the default "Normal" style so no stray style index is left on the cell.
$numRanges = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "D46", "E46", "D47", "D48", "E48", "D49", "E49", "D50", "D51")
foreach ($addr in $numRanges) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "330.19"
$ws.Range("E2").Value = "0.06%"
$ws.Range("D3").Value = "41.60"
$ws.Range("E3").Value = "1.19%"
$ws.Range("D4").Value = "5.680"
$ws.Range("E4").Value = "-0.35%"
$ws.Range("D5").Value = "0.08431"
$ws.Range("E5").Value = "4.49%"
$ws.Range("D6").Value = "8.797"
$ws.Range("E6").Value = "0.86%"
$ws.Range("D7").Value = "1.989"
$ws.Range("E7").Value = "-1.56%"
$ws.Range("D8").Value = "4.484"
$ws.Range("E8").Value = "-0.71%"
$ws.Range("D9").Value = "2.950"
$ws.Range("E9").Value = "0.51%"
$ws.Range("E10").Value = "0.43%"
$ws.Range("D11").Value = "0.1278"
$ws.Range("E11").Value = "0.38%"
$ws.Range("D12").Value = "0.1966"
$ws.Range("E12").Value = "1.20%"
$ws.Range("D13").Value = "0.09352"
$ws.Range("E13").Value = "-0.24%"
$ws.Range("D14").Value = "0.03955"
$ws.Range("E14").Value = "6.84%"
$ws.Range("D15").Value = "0.1062"
$ws.Range("E15").Value = "1.00%"
$ws.Range("D16").Value = "0.001315"
$ws.Range("E16").Value = "1.55%"
$ws.Range("D17").Value = "0.006112"
$ws.Range("E17").Value = "-2.03%"
$ws.Range("D18").Value = "0.004402"
$ws.Range("E18").Value = "0.19%"
$ws.Range("D19").Value = "3.423"
$ws.Range("E19").Value = "1.74%"
$ws.Range("D20").Value = "0.3500"
$ws.Range("E20").Value = "0.72%"
$ws.Range("D21").Value = "8.948"
$ws.Range("E21").Value = "8.55%"
$ws.Range("D22").Value = "0.1364"
$ws.Range("E22").Value = "-3.79%"
$ws.Range("D23").Value = "0.2513"
$ws.Range("E23").Value = "-5.20%"
$ws.Range("D24").Value = "0.04422"
$ws.Range("E24").Value = "-0.13%"
$ws.Range("D25").Value = "0.001247"
$ws.Range("E25").Value = "-1.13%"
$ws.Range("E26").Value = "-3.89%"
$ws.Range("D27").Value = "0.0003994"
$ws.Range("E27").Value = "0.04%"
$ws.Range("D39").Value = "0.02825"
$ws.Range("E39").Value = "-0.70%"
$ws.Range("D40").Value = "0.05513"
$ws.Range("E40").Value = "0.93%"
$ws.Range("D41").Value = "0.007904"
$ws.Range("E41").Value = "3.70%"
$ws.Range("D42").Value = "0.1437"
$ws.Range("E42").Value = "1.22%"
$ws.Range("D43").Value = "0.008989"
$ws.Range("E43").Value = "-9.61%"
$ws.Range("D44").Value = "0.002084"
$ws.Range("E44").Value = "-2.20%"
$ws.Range("E45").Value = "-7.33%"
$ws.Range("D46").Value = "0.00007227"
$ws.Range("E46").Value = "6.78%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").Value = "0.003260"
$ws.Range("E48").Value = "8.69%"
$ws.Range("D49").Value = "0.002281"
$ws.Range("E49").Value = "0.04%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("D51").Value = "0.0002003"

foreach ($addr in $numRanges) {
    $ws.Range($addr).Style = "Normal"
}
